# Alumni seminars workbook update:
#  - row156: Jazz (Freedom Rider) entry reclassified from "Upcoming" to "Hobby";
#            StartTime cleared; Slide/Video/Audio links added.
#  - row157: Recreational Auto Racing entry reclassified from "Upcoming" to
#            "Hobby"; StartTime cleared.
#  - three new rows inserted before the old last row (memorial tribute,
#    retirement-travel talk, old-village talk) and one new row appended at
#    the end (indie game talk).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# row156: "Hobby" reclass, clear start time
# ---------------------------------------------------------------------
$ws.Range("E156").Value = "Hobby"
$ws.Range("I156").ClearContents()

# ---------------------------------------------------------------------
# row157: "Hobby" reclass, clear start time
# ---------------------------------------------------------------------
$ws.Range("E157").Value = "Hobby"
$ws.Range("I157").ClearContents()

# ---------------------------------------------------------------------
# Make room for three new rows just before the (old) last row 158, which
# pushes the existing "Astronomy 3" row down to row 161. Its date moves
# from 12/13/2024 to 1/24/2025 to keep the bi-weekly schedule in order.
# ---------------------------------------------------------------------
$ws.Rows("158:160").Insert()
$ws.Range("A161").Value = 45681

# row158: memorial tribute talk (new row, picks up date/percent format
# automatically from the row Insert pushed down from below)
$ws.Range("A158").Value = 45618
$ws.Range("B158").Value = "All"
$ws.Range("D158").Value = "img/memorial_flowers.jpg"
$ws.Range("E158").Value = "Experience"

# row159: retirement/travel-life talk
$ws.Range("A159").Value = 45639
$ws.Range("B159").Value = "陈焕"
$ws.Range("E159").Value = "Upcoming"
$ws.Range("I159").Value = 0.8125

# row160: old-village talk
$ws.Range("A160").Value = 45667
$ws.Range("B160").Value = "何晖"
$ws.Range("E160").Value = "Upcoming"
$ws.Range("I160").Value = 0.8125

# ---------------------------------------------------------------------
# row162: brand-new last row (indie game talk), appended after the row
# that the earlier Insert pushed down to 161. Formats are copied from
# that row's A/I cells (date + percent-time styling) since row162 is
# beyond the sheet's previous used range and has no inherited format.
# ---------------------------------------------------------------------
$ws.Range("A161").Copy()
$ws.Range("A162").PasteSpecial(-4122)
$ws.Range("I161").Copy()
$ws.Range("I162").PasteSpecial(-4122)

$ws.Range("A162").Value = 45695
$ws.Range("E162").Value = "Upcoming"
$ws.Range("I162").Value = 0.8125

# ---------------------------------------------------------------------
# Remaining text fields, entered in this order so new shared-string table
# entries land in the same sequence as the source edit.
# ---------------------------------------------------------------------
$ws.Range("C159").Value = "退休后的旅居生活， 从沿海到边陲， 从繁华到中国最大的无人区"
$ws.Range("C160").Value = "古寨静谧如画卷，寻访村民烟火间"
$ws.Range("B162").Value = "Nicole He （贺敏女儿）"
$ws.Range("C162").Value = "My Journey to become an indie game creative director"
$ws.Range("C158").Value = "回忆赵凯华老师和我们度过的岁月"
$ws.Range("F156").Value = "https://drive.google.com/file/d/1LmHG8xGkfzdHotG3otUN5AWNG35KFKi4/view?usp=sharing"
$ws.Range("G156").Value = "https://drive.google.com/file/d/1WJ50WB2SARj-GEoFfLJQCuFgqsnqXZJr/view?usp=sharing"
$ws.Range("H156").Value = "https://drive.google.com/file/d/1DN3BZNtkXseE0CxtWLzmePr_5Npwn_gQ/view?usp=sharing"
$ws.Range("D159").Value = "img/旅居生活.jpg"
$ws.Range("D160").Value = "img/古寨.jpg"
$ws.Range("D162").Value = "img/crush_house_2.jpg"

# ---------------------------------------------------------------------
# Update the view: scroll so row 135 is at top and select I163 (the cell
# just past the new last row), matching where editing left off.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 135
$win.ScrollColumn = 1
$ws.Range("I163").Select()
